$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''25.555.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.80%  '

# Row 3
$ws.Range("D3").Value = '''1.669.73'
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = '''0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '''236.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.08%  '

# Row 6
$ws.Range("E6").Value = '  +0.00%  '

# Row 7
$ws.Range("D7").Value = '''0.4728'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.16%  '

# Row 8
$ws.Range("E8").Value = '  +1.55%  '

# Row 9
$ws.Range("D9").Value = '''0.06165'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.56%  '

# Row 10
$ws.Range("D10").Value = '''1.667.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.73%  '

# Row 11
$ws.Range("D11").Value = '''0.06993'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.77%  '

# Row 12
$ws.Range("D12").Value = '''14.78'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.60%  '

# Row 13
$ws.Range("D13").Value = '''0.5849'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.62%  '

# Row 14
$ws.Range("D14").Value = '''4.362'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.11%  '

# Row 15
$ws.Range("D15").Value = '''75.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.26%  '

# Row 16
$ws.Range("E16").Value = '  +0.01%  '

# Row 17
$ws.Range("D17").Value = '''0.9998'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.08%  '

# Row 18
$ws.Range("D18").Value = '''25.547.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.75%  '

# Row 19
$ws.Range("D19").Value = '''0.000006727'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.64%  '

# Row 20
$ws.Range("D20").Value = '''11.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.54%  '

# Row 21
$ws.Range("D21").Value = '''1.882.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.55%  '

# Row 22
$ws.Range("D22").Value = '''4.437'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.29%  '

# Row 23
$ws.Range("D23").Value = '''8.765'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.71%  '

# Row 24
$ws.Range("D24").Value = '''5.223'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.22%  '

# Row 25
$ws.Range("D25").Value = '''136.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.94%  '

# Row 26
$ws.Range("E26").Value = '  +1.62%  '

# Row 27
$ws.Range("D27").Value = '''1.385'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.67%  '

# Row 28
$ws.Range("D28").Value = '''1.712'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.99%  '

# Row 29
$ws.Range("D29").Value = '''104.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.20%  '

# Row 30
$ws.Range("D30").Value = '''3.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.33%  '

# Row 31
$ws.Range("D31").Value = '''0.07827'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.60%  '

# Row 32
$ws.Range("D32").Value = '''3.629'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.28%  '

# Row 33
$ws.Range("B33").Value = 'Frax'
$ws.Range("C33").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D33").Value = '''0.9993'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.04302'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.70%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.625'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.70%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''0.9534'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.84%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.6053'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.01%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''0.9361'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +16.15%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.525'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.38%  '

# Row 40
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '''0.9999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.14%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '''0.01478'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.74%  '

# Row 42
$ws.Range("D42").Value = '''100.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.16%  '

# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''1.839'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.68%  '

# Row 44
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.3743'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.76%  '

# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '''4.907'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.46%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '''0.1114'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.44%  '

# Row 47
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '''6.193'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.27%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.05264'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '

# Row 49
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''29.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.65%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''7.479'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.11%  '

# Row 51
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
